$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-05-13 Tuesday"; new="2025-05-14 Wednesday"},
    @{old="533×9="; new="989×6="},
    @{old="910×7="; new="488×6="},
    @{old="472×7="; new="127×5="},
    @{old="881×8="; new="918×3="},
    @{old="836×7="; new="520×9="},
    @{old="328×2="; new="487×6="},
    @{old="650×2="; new="684×9="},
    @{old="949×9="; new="541×2="},
    @{old="215×2="; new="354×4="},
    @{old="769×7="; new="540×4="},
    @{old="534×3="; new="606×7="},
    @{old="554×3="; new="509×2="},
    @{old="590×6="; new="836×2="},
    @{old="149×2="; new="153×4="},
    @{old="506×5="; new="957×7="},
    @{old="881×3="; new="274×6="},
    @{old="778×6="; new="156×3="},
    @{old="669×2="; new="575×6="},
    @{old="969×8="; new="295×3="},
    @{old="532×9="; new="162×3="},
    @{old="817×3="; new="985×6="},
    @{old="483×5="; new="867×2="},
    @{old="137×8="; new="532×3="},
    @{old="302×5="; new="759×3="},
    @{old="666×8="; new="142×6="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
